$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 557
$ws1.Range("F8").Value = 69
$ws1.Range("F9").Value = 6797
$ws1.Range("F10").Value = 160
$ws1.Range("F12").Value = 145
$ws1.Range("F15").Value = 1095
$ws1.Range("F16").Value = 16189
$ws1.Range("F20").Value = 182
$ws1.Range("F22").Value = 11353
$ws1.Range("F23").Value = 8
$ws1.Range("F24").Value = 989
$ws1.Range("F25").Value = 4469
$ws1.Range("F26").Value = 314
$ws1.Range("F27").Value = 388

# Sheet "全部类型" (sheet4.xml) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 557
$ws4.Range("F9").Value = 69
$ws4.Range("F10").Value = 6797
$ws4.Range("F11").Value = 160
$ws4.Range("F13").Value = 145
$ws4.Range("F17").Value = 1095
$ws4.Range("F18").Value = 16189
$ws4.Range("F22").Value = 183
$ws4.Range("F26").Value = 11353
$ws4.Range("F27").Value = 8
$ws4.Range("F28").Value = 989
$ws4.Range("F29").Value = 4469
$ws4.Range("F30").Value = 314
$ws4.Range("F31").Value = 388
